# EB adding KNN test and outliers checks
# also add add NALL to the dtatset
#
# Adds a new "result" row (row 52) to every summary sheet, mirroring the
# existing layout (column A = feature/stat name styled like the header
# column, column B = the corresponding value for that sheet).

$wb = $excel.ActiveWorkbook

# sheet name -> value to place in column B for the new "result" row
$values = @{
    "data_type" = "int64"
    "max"       = 2
    "min"       = 0
    "NAN"       = 0
    "unique"    = 3
}

foreach ($name in @("data_type", "max", "min", "NAN", "unique")) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A52").Value = "result"
    $ws.Range("B52").Value = $values[$name]

    # Copy the formatting (style) of row 51's label cell onto the new
    # label cell so A52 keeps the same bold/centered/bordered style as
    # every other entry in column A.
    $ws.Range("A51").Copy()
    $ws.Range("A52").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
